# Apply the updated S-curve / comet-1 projection values to columns C, E, F, G
# (D, the PLF column, is unchanged; E = C / D, F = E * k, G = C * k where k is the
# constant CO2-per-MJ emission factor already present in the workbook).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25
$ws.Range("C25").Value = 1.605796365659359
$ws.Range("E25").Value = 2.418585580357664
$ws.Range("F25").Value = 0.1771200564062623
$ws.Range("G25").Value = 0.117597138249909

# Row 26
$ws.Range("C26").Value = 1.595170373197115
$ws.Range("E26").Value = 2.423189895590971
$ws.Range("F26").Value = 0.1774572438022588
$ws.Range("G26").Value = 0.1168189659166369

# Row 31
$ws.Range("C31").Value = 1.559211995598527
$ws.Range("E31").Value = 2.259060943174633
$ws.Range("F31").Value = 0.1654376032545038
$ws.Range("G31").Value = 0.1141856293416303

# Row 32
$ws.Range("C32").Value = 1.567714000147012
$ws.Range("E32").Value = 2.289265470940818
$ws.Range("F32").Value = 0.1676495686714481
$ws.Range("G32").Value = 0.1148082558624463

# Row 33
$ws.Range("C33").Value = 1.564079939769674
$ws.Range("E33").Value = 2.264544154068898
$ws.Range("F33").Value = 0.1658391547359842
$ws.Range("G33").Value = 0.1145421230514987

# Row 34
$ws.Range("C34").Value = 1.554941631335987
$ws.Range("E34").Value = 2.187896682526145
$ws.Range("F34").Value = 0.1602260374688904
$ws.Range("G34").Value = 0.1138728981465057

# Row 35
$ws.Range("C35").Value = 1.550489742549896
$ws.Range("E35").Value = 2.239272596246528
$ws.Range("F35").Value = 0.1639884450553657
$ws.Range("G35").Value = 0.1135468733825648

# Row 38
$ws.Range("C38").Value = 1.427985993580319
$ws.Range("E38").Value = 1.945719472026729
$ws.Range("F38").Value = 0.1424906959815635
$ws.Range("G38").Value = 0.1045755675484081

# Row 39
$ws.Range("C39").Value = 1.410548045123722
$ws.Range("E39").Value = 1.882167269668976
$ws.Range("F39").Value = 0.1378365833639389
$ws.Range("G39").Value = 0.1032985358653757

# Row 40
$ws.Range("C40").Value = 1.405378496021678
$ws.Range("E40").Value = 1.853692822352101
$ws.Range("F40").Value = 0.1357513167701655
$ws.Range("G40").Value = 0.1029199547492121

# Row 41
$ws.Range("C41").Value = 1.408893757585069
$ws.Range("E41").Value = 1.834384170021593
$ws.Range("F41").Value = 0.1343372880015813
$ws.Range("G41").Value = 0.1031773875774929

# Row 42
$ws.Range("C42").Value = 1.380351075228878
$ws.Range("E42").Value = 1.815920663940857
$ws.Range("F42").Value = 0.1329851517509411
$ws.Range("G42").Value = 0.1010871239333315

# Row 43
$ws.Range("C43").Value = 1.363977034209024
$ws.Range("E43").Value = 1.778751645324429
$ws.Range("F43").Value = 0.1302631564130984
$ws.Range("G43").Value = 0.09988800528622285

# Row 44
$ws.Range("C44").Value = 1.380901000651561
$ws.Range("E44").Value = 1.766531543861792
$ws.Range("F44").Value = 0.1293682428413271
$ws.Range("G44").Value = 0.1011273965714701

# Row 47
$ws.Range("C47").Value = 1.342218401785662
$ws.Range("E47").Value = 1.688706241602388
$ws.Range("F47").Value = 0.1236688696051807
$ws.Range("G47").Value = 0.0982945573497727

# Row 48
$ws.Range("C48").Value = 1.331246744782096
$ws.Range("E48").Value = 1.669934466960038
$ws.Range("F48").Value = 0.122294157951187
$ws.Range("G48").Value = 0.0974910709967885

# Row 49
$ws.Range("C49").Value = 1.315140810373198
$ws.Range("E49").Value = 1.639036470100827
$ws.Range("F49").Value = 0.1200314077756342
$ws.Range("G49").Value = 0.09631158657657721

# Row 50
$ws.Range("C50").Value = 1.304783119122762
$ws.Range("E50").Value = 1.625178676835295
$ws.Range("F50").Value = 0.1190165612699776
$ws.Range("G50").Value = 0.09555306272138883

# Row 51
$ws.Range("C51").Value = 1.302501009379456
$ws.Range("E51").Value = 1.598601300647098
$ws.Range("F51").Value = 0.1170702227125106
$ws.Range("G51").Value = 0.09538593718746423

# Row 52
$ws.Range("C52").Value = 1.296052864122837
$ws.Range("E52").Value = 1.588232197226381
$ws.Range("F52").Value = 0.1163108631109007
$ws.Range("G52").Value = 0.09491372075615678

# Row 53
$ws.Range("C53").Value = 1.276281382275522
$ws.Range("E53").Value = 1.549199684641708
$ws.Range("F53").Value = 0.1134523986898679
$ws.Range("G53").Value = 0.09346579763593625

# Row 54
$ws.Range("C54").Value = 1.190921954421548
$ws.Range("E54").Value = 1.825434111347368
$ws.Range("F54").Value = 0.1336818491739903
$ws.Range("G54").Value = 0.08721467846980514

# Row 55
$ws.Range("C55").Value = 1.150148931499757
$ws.Range("E55").Value = 1.699512378614168
$ws.Range("F55").Value = 0.1244602344477583
$ws.Range("G55").Value = 0.08422875141458244

# Row 56
$ws.Range("C56").Value = 1.142005350412699
$ws.Range("E56").Value = 1.681497877118
$ws.Range("F56").Value = 0.1231409801087573
$ws.Range("G56").Value = 0.08363237328630661

# Row 57
$ws.Range("C57").Value = 1.135647554608852
$ws.Range("E57").Value = 1.667433730055871
$ws.Range("F57").Value = 0.1221110217144045
$ws.Range("G57").Value = 0.08316677340820329

# Row 58
$ws.Range("C58").Value = 1.124752312390587
$ws.Range("E58").Value = 1.643332249989076
$ws.Range("F58").Value = 0.1203460002309498
$ws.Range("G58").Value = 0.08236888313219592

# Row 59
$ws.Range("C59").Value = 1.106894511493559
$ws.Range("E59").Value = 1.603828820534797
$ws.Range("F59").Value = 0.1174530491979133
$ws.Range("G59").Value = 0.08106110443382726

# Row 60
$ws.Range("C60").Value = 1.079657624613052
$ws.Range("E60").Value = 1.543577816189835
$ws.Range("F60").Value = 0.1130406929121641
$ws.Range("G60").Value = 0.07906646799020264

# Row 61
$ws.Range("C61").Value = 1.042361339290171
$ws.Range("E61").Value = 1.461074310598378
$ws.Range("F61").Value = 0.1069987212396495
$ws.Range("G61").Value = 0.07633515254129643

# Row 62
$ws.Range("C62").Value = 0.9981777361154145
$ws.Range("E62").Value = 1.363335296645078
$ws.Range("F62").Value = 0.09984100897794779
$ws.Range("G62").Value = 0.07309945877461668

# Row 63
$ws.Range("C63").Value = 0.9539941329406577
$ws.Range("E63").Value = 1.265596282691778
$ws.Range("F63").Value = 0.09268329671624609
$ws.Range("G63").Value = 0.06986376500793694

# Row 64
$ws.Range("C64").Value = 0.9166978476177772
$ws.Range("E64").Value = 1.18309277710032
$ws.Range("F64").Value = 0.08664132504373147
$ws.Range("G64").Value = 0.06713244955903074

# Row 65
$ws.Range("C65").Value = 0.8894609607372701
$ws.Range("E65").Value = 1.122841772755358
$ws.Range("F65").Value = 0.08222896875798222
$ws.Range("G65").Value = 0.06513781311540612

# Row 66
$ws.Range("C66").Value = 0.8716031598402422
$ws.Range("E66").Value = 1.083338343301079
$ws.Range("F66").Value = 0.07933601772494581
$ws.Range("G66").Value = 0.06383003441703744

# Row 67
$ws.Range("C67").Value = 0.8607079176219773
$ws.Range("E67").Value = 1.059236863234285
$ws.Range("F67").Value = 0.07757099624149107
$ws.Range("G67").Value = 0.06303214414103009

# Row 68
$ws.Range("C68").Value = 0.8543501218181296
$ws.Range("E68").Value = 1.045172716172156
$ws.Range("F68").Value = 0.07654103784713821
$ws.Range("G68").Value = 0.06256654426292677

# Row 69
$ws.Range("C69").Value = 0.8507361338896253
$ws.Range("E69").Value = 1.037178175410642
$ws.Range("F69").Value = 0.07595557437537961
$ws.Range("G69").Value = 0.06230188141578717

# Row 70
$ws.Range("C70").Value = 0.8462065407310719
$ws.Range("E70").Value = 1.027158214675988
$ws.Range("F70").Value = 0.07522178350813724
$ws.Range("G70").Value = 0.06197016613465093

# Row 71
$ws.Range("C71").Value = 0.8462065407310719
$ws.Range("E71").Value = 1.027158214675988
$ws.Range("F71").Value = 0.07522178350813724
$ws.Range("G71").Value = 0.06197016613465093

# Row 72
$ws.Range("C72").Value = 0.8462065407310719
$ws.Range("E72").Value = 1.027158214675988
$ws.Range("F72").Value = 0.07522178350813724
$ws.Range("G72").Value = 0.06197016613465093

# Row 73
$ws.Range("C73").Value = 0.8462065407310719
$ws.Range("E73").Value = 1.027158214675988
$ws.Range("F73").Value = 0.07522178350813724
$ws.Range("G73").Value = 0.06197016613465093

# Row 74
$ws.Range("C74").Value = 0.8452240531960947
$ws.Range("E74").Value = 1.025965633321681
$ws.Range("F74").Value = 0.07513444730698753
$ws.Range("G74").Value = 0.06189821571493997

# Row 75
$ws.Range("C75").Value = 0.8441160234904356
$ws.Range("E75").Value = 1.024620664026963
$ws.Range("F75").Value = 0.07503595129374746
$ws.Range("G75").Value = 0.06181707147693573

# Row 76
$ws.Range("C76").Value = 0.8418709193803046
$ws.Range("E76").Value = 1.021895469859199
$ws.Range("F76").Value = 0.07483637739872696
$ws.Range("G76").Value = 0.06165265597315788

# Row 77
$ws.Range("C77").Value = 0.8376531692092239
$ws.Range("E77").Value = 1.01677580163738
$ws.Range("F77").Value = 0.07446144920449872
$ws.Range("G77").Value = 0.06134377786097678

# Row 78
$ws.Range("C78").Value = 0.8307484465687006
$ws.Range("E78").Value = 1.008394582350011
$ws.Range("F78").Value = 0.07384766813965318
$ws.Range("G78").Value = 0.06083812493990948

# Row 79
$ws.Range("C79").Value = 0.8216443523566421
$ws.Range("E79").Value = 0.9973436808183105
$ws.Range("F79").Value = 0.07303837847939423
$ws.Range("G79").Value = 0.06017140564187692

# Row 80
$ws.Range("C80").Value = 0.8125402581445836
$ws.Range("E80").Value = 0.98629277928661
$ws.Range("F80").Value = 0.07222908881913528
$ws.Range("G80").Value = 0.05950468634384437

# Row 81
$ws.Range("C81").Value = 0.8056355355040603
$ws.Range("E81").Value = 0.9779115599992413
$ws.Range("F81").Value = 0.07161530775428974
$ws.Range("G81").Value = 0.05899903342277708

# Row 82
$ws.Range("C82").Value = 0.8014177853329798
$ws.Range("E82").Value = 0.9727918917774221
$ws.Range("F82").Value = 0.0712403795600615
$ws.Range("G82").Value = 0.05869015531059597

# Row 83
$ws.Range("C83").Value = 0.7991726812228487
$ws.Range("E83").Value = 0.9700666976096581
$ws.Range("F83").Value = 0.07104080566504102
$ws.Range("G83").Value = 0.05852573980681813

# Row 84
$ws.Range("C84").Value = 0.7980646515171895
$ws.Range("E84").Value = 0.9687217283149402
$ws.Range("F84").Value = 0.07094230965180096
$ws.Range("G84").Value = 0.05844459556881389

